# Auto commit - 08121738
# Appends 4 new maintenance rows (54-57) to the 'Report' sheet, mirroring
# the existing table's alternating row styling, and extends the print
# area / dimension / selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Clone formatting (fill/border/alignment) from the two alternating
#    "template" rows before we touch row 53's own formatting, so the new
#    row 57 (which must keep the *current*, non-wrapped look of row 53)
#    is copied first.
# ---------------------------------------------------------------------
$ws.Range('A53:AK53').Copy($ws.Range('A57:AK57'))   # odd/shaded template, no wrap (matches target row57)
$ws.Range('A52:AK52').Copy($ws.Range('A54:AK54'))   # even/plain template
$ws.Range('A53:AK53').Copy($ws.Range('A55:AK55'))   # odd/shaded template
$ws.Range('A52:AK52').Copy($ws.Range('A56:AK56'))   # even/plain template

# ---------------------------------------------------------------------
# 2) Fix up the long-standing wrap-text anomaly on row 53 (P/AC were not
#    wrapping) and apply the same wrap formatting to the new interior
#    rows (54-56). Row 57, the new last row, intentionally keeps the
#    no-wrap look that row 53 used to have.
# ---------------------------------------------------------------------
$ws.Range('P53').WrapText = $true
$ws.Range('AC53').WrapText = $true
$ws.Range('P54').WrapText = $true
$ws.Range('AC54').WrapText = $true
$ws.Range('P55').WrapText = $true
$ws.Range('AC55').WrapText = $true
$ws.Range('P56').WrapText = $true
$ws.Range('AC56').WrapText = $true

# ---------------------------------------------------------------------
# 3) Row 54 - D620 / 北縣三愛三店
# ---------------------------------------------------------------------
$ws.Range('A54').Value = 52
$ws.Range('C54').Value = 2025081485
$ws.Range('F54').Value = 'D620'
$ws.Range('G54').Value = '北縣三愛三店'
$ws.Range('H54').Value = '新北市三重區'
$ws.Range('Q54').Value = 'THILF0D620'
$ws.Range('S54').Value = '吳宗鴻'
$ws.Range('V54').Value = '2025-08-12 14:19:46'
$ws.Range('W54').Value = '2025-08-12 14:00:00'
$ws.Range('X54').Value = '2025-08-12 14:19:00'
$ws.Range('Z54').Value = 0.3

# ---------------------------------------------------------------------
# 4) Row 55 - 3796 / 三重進安店
# ---------------------------------------------------------------------
$ws.Range('A55').Value = 53
$ws.Range('C55').Value = 2025081493
$ws.Range('F55').Value = 3796
$ws.Range('G55').Value = '三重進安店'
$ws.Range('H55').Value = '新北市三重區'
$ws.Range('Q55').Value = 'THILF03796'
$ws.Range('S55').Value = '吳宗鴻'
$ws.Range('V55').Value = '2025-08-12 14:48:58'
$ws.Range('W55').Value = '2025-08-12 14:30:00'
$ws.Range('X55').Value = '2025-08-12 14:47:00'
$ws.Range('Z55').Value = 0.3

# ---------------------------------------------------------------------
# 5) Row 56 - 3811 / 三重碧華公園
# ---------------------------------------------------------------------
$ws.Range('A56').Value = 54
$ws.Range('C56').Value = 2025081505
$ws.Range('F56').Value = 3811
$ws.Range('G56').Value = '三重碧華公園'
$ws.Range('H56').Value = '新北市三重區'
$ws.Range('Q56').Value = 'THILF03811'
$ws.Range('S56').Value = '吳宗鴻'
$ws.Range('V56').Value = '2025-08-12 15:29:11'
$ws.Range('W56').Value = '2025-08-12 15:15:00'
$ws.Range('X56').Value = '2025-08-12 15:30:00'
$ws.Range('Z56').Value = 0.3
$ws.Range('AC56').Value = 'PMQ3+TVV'

# ---------------------------------------------------------------------
# 6) Row 57 - D350 / 三重徐匯店
# ---------------------------------------------------------------------
$ws.Range('A57').Value = 55
$ws.Range('C57').Value = 2025081518
$ws.Range('F57').Value = 'D350'
$ws.Range('G57').Value = '三重徐匯店'
$ws.Range('H57').Value = '新北市三重區'
$ws.Range('Q57').Value = 'THILF0D350'
$ws.Range('S57').Value = '吳宗鴻'
$ws.Range('V57').Value = '2025-08-12 16:23:19'
$ws.Range('W57').Value = '2025-08-12 16:00:00'
$ws.Range('X57').Value = '2025-08-12 16:22:00'
$ws.Range('Z57').Value = 0.4

# ---------------------------------------------------------------------
# 7) Extend the print area to the new last row and move the selection
#    to match Excel's usual "active cell follows last edit" behaviour.
# ---------------------------------------------------------------------
$wb.Names.Item('Report!Print_Area').RefersTo = "='Report'!`$A`$1:`$AK`$57"
$ws.Range('A57').Select()
